$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Transfer Volume for the A2 destination well row (row 3) to 2875
$ws.Range("H3").Value = 2875

# Remove rows 4-7 (UID 3,4,5,6 / destination wells A3-A6) - no longer needed
$ws.Range("A4:I7").EntireRow.Delete()
